$wb = $excel.ActiveWorkbook

# --- Rename sheets --------------------------------------------------------
$wsAufwand = $wb.Worksheets.Item(1)
$wsAufwand.Name = "Aufwand"

$wsRessourcen = $wb.Worksheets.Item(2)
$wsRessourcen.Name = "Ressourcen"

# --- Remove the unused third sheet ----------------------------------------
$wb.Worksheets.Item(3).Delete()

# --- Aufwand sheet tweaks ---------------------------------------------------
$wsAufwand.Range("A163").Value = "TOTAL Aufwand"

# --- Fill in the new Ressourcen sheet --------------------------------------
$ws = $wsRessourcen

$ws.Range("A1").Value = 2014

$ws.Range("B2").Value = "Html UI %"
$ws.Range("C2").Value = "Ferien (Wo.)"
$ws.Range("D2").Value = "Html UI PT"

$ws.Range("A3").Value = "AWE"
$ws.Range("B3").Value = 80
$ws.Range("C3").Value = 9
$ws.Range("D3").Formula = '=$B$10*0.01*B3-(C3*5)'

$ws.Range("A4").Value = "BSH"
$ws.Range("B4").Value = 90
$ws.Range("C4").Value = 6
$ws.Range("D4").Formula = '=$B$10*0.01*B4-(C4*5)'

$ws.Range("A5").Value = "CGU"
$ws.Range("B5").Value = 85
$ws.Range("C5").Value = 5
$ws.Range("D5").Formula = '=$B$10*0.01*B5-(C5*5)'

$ws.Range("A6").Value = "CRU"
$ws.Range("B6").Value = 40
$ws.Range("C6").Value = 2
$ws.Range("D6").Formula = '=$B$10*0.01*B6-(C6*5)'

$ws.Range("A7").Value = "IMO"
$ws.Range("B7").Value = 40
$ws.Range("C7").Value = 5
$ws.Range("D7").Formula = '=$B$10*0.01*B7-(C7*5)'

$ws.Range("A8").Value = "Total"
$ws.Range("D8").Formula = "=SUM(D3:D7)"

$ws.Range("A10").Value = "Arbeitstage"
$ws.Range("B10").Value = 150
$ws.Range("C10").Value = "bis Ende Jahr"

$ws.Range("A12").Value = "2015 (prov.)"

$ws.Range("B13").Value = "Html UI %"
$ws.Range("C13").Value = "Ferien (Wo.)"
$ws.Range("D13").Value = "Html UI PT"

$ws.Range("A14").Value = "AWE"
$ws.Range("B14").Value = 80
$ws.Range("C14").Value = 6
$ws.Range("D14").Formula = '=$B$21*0.01*B14-(C14*5)'

$ws.Range("A15").Value = "BSH"
$ws.Range("B15").Value = 90
$ws.Range("C15").Value = 6
$ws.Range("D15").Formula = '=$B$21*0.01*B15-(C15*5)'

$ws.Range("A16").Value = "CGU"
$ws.Range("B16").Value = 85
$ws.Range("C16").Value = 6
$ws.Range("D16").Formula = '=$B$21*0.01*B16-(C16*5)'

$ws.Range("A17").Value = "CRU"
$ws.Range("B17").Value = 40
$ws.Range("C17").Value = 6
$ws.Range("D17").Formula = '=$B$21*0.01*B17-(C17*5)'

$ws.Range("A18").Value = "IMO"
$ws.Range("B18").Value = 40
$ws.Range("C18").Value = 6
$ws.Range("D18").Formula = '=$B$21*0.01*B18-(C18*5)'

$ws.Range("A19").Value = "Total"
$ws.Range("D19").Formula = "=SUM(D14:D18)"

$ws.Range("A21").Value = "Arbeitstage"
$ws.Range("B21").Value = 240
$ws.Range("C21").Value = "bis Ende Jahr"

$ws.Range("A23").Value = "Prognose: Umsetzung gemäss Schätzung (826 PT) Ende 3. Quartal 2015"

# --- Activate the Ressourcen tab, set its selection -------------------------
$wsRessourcen.Activate()
$wsRessourcen.Range("G17").Select()
